$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Remove the rows for festivals that were dropped from the sheet.
# Deleted from the bottom up so earlier row numbers stay valid.
# Row numbers below refer to the ORIGINAL (pre-edit) layout:
#   35 TCU Amphitheater at White River State Park
#   24 IndyFringe Theatre
#   21 Indianapolis Motor Speedway
#   18 Indiana State Fairgrounds & Event Center
#   16 Indiana Latino Expo
#   13 Hogan Farms Pumpkin Patch & Corn Maze
#    8 Conner Prairie
# ------------------------------------------------------------------
$ws.Rows.Item(35).Delete()
$ws.Rows.Item(24).Delete()
$ws.Rows.Item(21).Delete()
$ws.Rows.Item(18).Delete()
$ws.Rows.Item(16).Delete()
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(8).Delete()

# ------------------------------------------------------------------
# Insert a new row for "Free Stage", which now sits between
# "Festival Flea Market" (row 9) and "HSI Show Productions" (row 11).
# ------------------------------------------------------------------
$ws.Rows.Item(10).Insert()

$rankCell = $ws.Cells.Item(10, 1)
$rankCell.Value = 26
$rankCell.Borders.LineStyle = 1
$rankCell.HorizontalAlignment = -4108
$rankCell.VerticalAlignment = -4160
$rankCell.Font.Bold = $true

$ws.Cells.Item(10, 2).Value = "OPERATIONAL"
$ws.Cells.Item(10, 3).Value = "Free Stage"
$ws.Cells.Item(10, 4).Value = 4.6
$ws.Cells.Item(10, 5).Value = 361

# ------------------------------------------------------------------
# Refresh the rank values in column A for every remaining row, since
# the source ranking was recomputed after the festivals above were
# removed from the dataset.
# ------------------------------------------------------------------
$ws.Cells.Item(2, 1).Value = 10
$ws.Cells.Item(3, 1).Value = 11
$ws.Cells.Item(4, 1).Value = 6
$ws.Cells.Item(5, 1).Value = 28
$ws.Cells.Item(6, 1).Value = 8
$ws.Cells.Item(7, 1).Value = 23
$ws.Cells.Item(8, 1).Value = 25
$ws.Cells.Item(9, 1).Value = 24
$ws.Cells.Item(10, 1).Value = 26
$ws.Cells.Item(11, 1).Value = 17
$ws.Cells.Item(12, 1).Value = 22
$ws.Cells.Item(13, 1).Value = 27
$ws.Cells.Item(14, 1).Value = 21
$ws.Cells.Item(15, 1).Value = 9
$ws.Cells.Item(16, 1).Value = 16
$ws.Cells.Item(17, 1).Value = 18
$ws.Cells.Item(18, 1).Value = 29
$ws.Cells.Item(19, 1).Value = 0
$ws.Cells.Item(20, 1).Value = 1
$ws.Cells.Item(21, 1).Value = 32
$ws.Cells.Item(22, 1).Value = 15
$ws.Cells.Item(23, 1).Value = 19
$ws.Cells.Item(24, 1).Value = 13
$ws.Cells.Item(25, 1).Value = 30
$ws.Cells.Item(26, 1).Value = 14
$ws.Cells.Item(27, 1).Value = 12
$ws.Cells.Item(28, 1).Value = 3
$ws.Cells.Item(29, 1).Value = 4
$ws.Cells.Item(30, 1).Value = 5
$ws.Cells.Item(31, 1).Value = 31
$ws.Cells.Item(32, 1).Value = 7
$ws.Cells.Item(33, 1).Value = 20
$ws.Cells.Item(34, 1).Value = 2

# ------------------------------------------------------------------
# A couple of underlying data points were also refreshed:
#   Indianapolis Zoo: user_ratings_total 14979 -> 14982
#   Waterman's Family Farm: rating 4.3 -> 4.4, user_ratings_total 668 -> 669
# ------------------------------------------------------------------
$ws.Cells.Item(18, 5).Value = 14982
$ws.Cells.Item(33, 4).Value = 4.4
$ws.Cells.Item(33, 5).Value = 669
